$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I6").Value = 1.477282346912468
$ws.Range("J6").Value = 0.6430008601316372
$ws.Range("K6").Value = -0.2172062532435262
$ws.Range("L6").Value = 2.742617893978283
